$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187 (shifts the existing rows 187-213 down to 188-214)
$ws.Rows(187).Insert()

# Populate the newly inserted row 187 with the new price-report entry
$ws.Cells.Item(187, 1).Value  = 1
$ws.Cells.Item(187, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(187, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(187, 4).Value  = 45244
$ws.Cells.Item(187, 5).Value  = 15
$ws.Cells.Item(187, 6).Value  = "Fruta"
$ws.Cells.Item(187, 7).Value  = 100102
$ws.Cells.Item(187, 8).Value  = "Cítricos"
$ws.Cells.Item(187, 9).Value  = 100102004
$ws.Cells.Item(187, 10).Value = "Mandarina"
$ws.Cells.Item(187, 11).Value = "Murcott"
$ws.Cells.Item(187, 12).Value = "Segunda"
$ws.Cells.Item(187, 13).Value = 270
$ws.Cells.Item(187, 14).Value = 15000
$ws.Cells.Item(187, 15).Value = 17000
$ws.Cells.Item(187, 16).Value = 16000
$ws.Cells.Item(187, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(187, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(187, 19).Value = 800
$ws.Cells.Item(187, 20).Value = 20
